$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert two new worksheets ("line1" and "Line2") between "Lines" and "Load"
# ---------------------------------------------------------------------------
$lines = $wb.Worksheets.Item("Lines")

$line1 = $wb.Worksheets.Add($lines.Next)
$line1.Name = "line1"

$line2 = $wb.Worksheets.Add($line1.Next)
$line2.Name = "Line2"

# ---------------------------------------------------------------------------
# 2. Populate "line1" - header row + the line that used to be row 2 of "Lines"
#    (From bus 1, To bus 2)
# ---------------------------------------------------------------------------
$line1.Range("B1").Value = "From bus"
$line1.Range("C1").Value = "To bus"
$line1.Range("D1").Value = "Length_km"
$line1.Range("E1").Value = "r_ohm_per_km"
$line1.Range("F1").Value = "x_ohm_per_km"
$line1.Range("G1").Value = "c_nf_per_km"
$line1.Range("H1").Value = "Max_current_kA"

$line1.Range("A2").Value = 0
$line1.Range("B2").Value = 1
$line1.Range("C2").Value = 2
$line1.Range("D2").Value = 125
$line1.Range("E2").Value = 0.6
$line1.Range("F2").Value = 0.08
$line1.Range("G2").Value = 210
$line1.Range("H2").Value = 0.142

$line1.Columns.Item(4).ColumnWidth = 14.166666666666668
$line1.Columns.Item(5).ColumnWidth = 13
$line1.Columns.Item(6).ColumnWidth = 11.666666666666668
$line1.Columns.Item(7).ColumnWidth = 11.333333333333332
$line1.Columns.Item(8).ColumnWidth = 14.5

$line1.Activate() | Out-Null
$line1.Range("E8").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3. Populate "Line2" - header row + the two lines that used to be rows 3 & 4
#    of "Lines" (From bus 3, To bus 4) and (From bus 4, To bus 2)
# ---------------------------------------------------------------------------
$line2.Range("B1").Value = "From bus"
$line2.Range("C1").Value = "To bus"
$line2.Range("D1").Value = "Length_km"
$line2.Range("E1").Value = "r_ohm_per_km"
$line2.Range("F1").Value = "x_ohm_per_km"
$line2.Range("G1").Value = "c_nf_per_km"
$line2.Range("H1").Value = "Max_current_kA"

$line2.Range("B2").Value = 3
$line2.Range("C2").Value = 4
$line2.Range("D2").Value = 60
$line2.Range("E2").Value = 0.6
$line2.Range("F2").Value = 0.08
$line2.Range("G2").Value = 210
$line2.Range("H2").Value = 0.142

$line2.Range("B3").Value = 4
$line2.Range("C3").Value = 2
$line2.Range("D3").Value = 65
$line2.Range("E3").Value = 0.6
$line2.Range("F3").Value = 0.08
$line2.Range("G3").Value = 210
$line2.Range("H3").Value = 0.142

$line2.Columns.Item(4).ColumnWidth = 13.333333333333332
$line2.Columns.Item(5).ColumnWidth = 15.5
$line2.Columns.Item(6).ColumnWidth = 17.5
$line2.Columns.Item(7).ColumnWidth = 13.166666666666668
$line2.Columns.Item(8).ColumnWidth = 15.666666666666668
$line2.Columns.Item(9).ColumnWidth = 13.666666666666668

$line2.Activate() | Out-Null
$line2.Range("D1:D1048576").Select() | Out-Null

# ---------------------------------------------------------------------------
# 4. Update "Lines" sheet data (renumbered bus references) and append the
#    divider row 30.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Lines")

$ws.Range("C2").Value = 2
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 4
$ws.Range("C4").Value = 2
$ws.Range("B5").Value = 6
$ws.Range("C5").Value = 7
$ws.Range("B6").Value = 8
$ws.Range("C6").Value = 2

$ws.Range("D30").Value = "------------------------------------------------------"

$ws.Activate() | Out-Null
$ws.Range("B1:H1").Select() | Out-Null

# ---------------------------------------------------------------------------
# 5. Update "Load" sheet (bus reference renumbered)
# ---------------------------------------------------------------------------
$load = $wb.Worksheets.Item("Load")
$load.Range("B2").Value = 2

$load.Activate()
$load.Range("N18").Select()

# ---------------------------------------------------------------------------
# 6. Leave "Line2" as the active sheet/tab (matches target activeTab)
# ---------------------------------------------------------------------------
$line2.Activate()
